# Updated symbol list on Thu Jan 12 13:55:40 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) readings for
# the existing coin rows, and patches the handful of rows whose Coin/Link
# (columns B/C) shifted because the upstream ranking re-sorted that day.
#
# Numeric-looking values are entered with a leading apostrophe so Excel
# keeps them as literal text (matching "283.53", "1.90%", etc. being
# plain strings in the sheet) instead of auto-converting them to numbers
# or percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'284.74"
$ws.Range("E2").Value = "'2.36%"

# Row 3 - OKB
$ws.Range("D3").Value = "'28.32"
$ws.Range("E3").Value = "'3.97%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.020"
$ws.Range("E4").Value = "'3.40%"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.06538"
$ws.Range("E5").Value = "'2.07%"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "'7.249"
$ws.Range("E6").Value = "'3.72%"

# Row 7 - FTXToken
$ws.Range("E7").Value = "'19.15%"

# Row 8 - MXToken
$ws.Range("E8").Value = "'5.01%"

# Row 9 - WazirX
$ws.Range("D9").Value = "'0.1552"
$ws.Range("E9").Value = "'1.86%"

# Row 10 - LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.06479"
$ws.Range("E10").Value = "'25.35%"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value = "'0.07563"
$ws.Range("E11").Value = "'0.77%"

# Row 12 - BitrueCoin
$ws.Range("D12").Value = "'0.02761"
$ws.Range("E12").Value = "'-5.82%"

# Row 13 - BitMartToken
$ws.Range("D13").Value = "'0.08947"
$ws.Range("E13").Value = "'-0.27%"

# Row 14 - BitForexToken
$ws.Range("D14").Value = "'0.001599"
$ws.Range("E14").Value = "'1.97%"

# Row 15 - was "One", now "CoinExToken" (ranking shifted up one slot)
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04400"
$ws.Range("E15").Value = "'-0.12%"

# Row 16 - was "TigerCash", now "One"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0006361"
$ws.Range("E16").Value = "'-0.22%"

# Row 17 - was "LEO", now "TigerCash"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005980"
$ws.Range("E17").Value = "'-1.65%"

# Row 18 - was "GateToken", now "LEO"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.443"
$ws.Range("E18").Value = "'-0.89%"

# Row 19 - was "BTSEToken", now "GateToken"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.375"
$ws.Range("E19").Value = "'2.17%"

# Row 20 - was "BitpandaEcosystemToken", now "BTSEToken"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.235"
$ws.Range("E20").Value = "'-0.47%"

# Row 21 - was "ProBitToken", now "BitpandaEcosystemToken"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "'0.3194"
$ws.Range("E21").Value = "'3.55%"

# Row 22 - was "MCDex", now "ProBitToken"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1306"
$ws.Range("E22").Value = "'-1.44%"

# Row 23 - was "ZBToken", now "MCDex"
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D23").Value = "'3.977"
$ws.Range("E23").Value = "'1.38%"

# Row 24 - was "CoinExToken", now "ZBToken"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "'0.1539"
$ws.Range("E24").Value = "'2.19%"

# Row 25 - BitKan
$ws.Range("D25").Value = "'0.001180"
$ws.Range("E25").Value = "'0.38%"

# Row 26 - HotbitToken
$ws.Range("D26").Value = "'0.004430"
$ws.Range("E26").Value = "'13.71%"

# Row 27 - NitroEx
$ws.Range("D27").Value = "'0.0001247"
$ws.Range("E27").Value = "'5.70%"

# Row 28 - UpBots
$ws.Range("D28").Value = "'0.0001615"
$ws.Range("E28").Value = "'-1.80%"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.04153"
$ws.Range("E40").Value = "'1.82%"

# Row 41 - KickToken
$ws.Range("D41").Value = "'0.006675"
$ws.Range("E41").Value = "'-1.97%"

# Row 42 - BKEXToken
$ws.Range("D42").Value = "'0.1228"
$ws.Range("E42").Value = "'4.63%"

# Row 43 - CEJI
$ws.Range("D43").Value = "'0.002164"
$ws.Range("E43").Value = "'14.54%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.01144"
$ws.Range("E44").Value = "'2.08%"

# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00005604"
$ws.Range("E45").Value = "'4.56%"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "'0.01297"
$ws.Range("E47").Value = "'-30.00%"
